$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Q1 date value (Sep 4, 2024 -> Oct 4, 2024)
$ws.Range("Q1").Value = 45569

# Add new course text "Fonts in css" into Q2, copying the formatting
# (the "Good" style + border) from an existing cell that already has it
$ws.Range("C2").Copy()
$ws.Range("Q2").PasteSpecial(-4122)
$ws.Range("Q2").Value = "Fonts in css"
$excel.CutCopyMode = 0

# Update the view: selecting Q2 also resets the scrolled topLeftCell
$ws.Range("Q2").Select()
